$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index map: B=2, C=3, D=4, E=5, F=6, G=7
$colIndex = @{ B=2; C=3; D=4; E=5; G=7 }

# Per-row updates parsed from the coinranking.com symbol-list refresh.
# Columns D, E and G hold numeric-looking text (price / % change / hour) that
# must stay stored as text, exactly like the source sheet (inline strings),
# so we force the Text number format before writing those values.
$updates = @(
    @{Row=2; D="257.03"; E="-0.75%"; G="19"},
    @{Row=3; D="27.36"; E="-2.05%"; G="19"},
    @{Row=4; D="4.585"; E="-11.95%"; G="19"},
    @{Row=5; D="0.05888"; E="-0.89%"; G="19"},
    @{Row=6; E="-1.17%"; G="19"},
    @{Row=7; D="0.8568"; E="-1.76%"; G="19"},
    @{Row=8; D="0.9234"; E="-7.36%"; G="19"},
    @{Row=9; D="0.1410"; E="-1.26%"; G="19"},
    @{Row=10; D="0.03582"; E="-1.04%"; G="19"},
    @{Row=11; D="0.07078"; E="-2.58%"; G="19"},
    @{Row=12; D="0.03213"; E="-1.21%"; G="19"},
    @{Row=13; D="0.09205"; E="-0.26%"; G="19"},
    @{Row=14; D="0.001553"; E="-1.22%"; G="19"},
    @{Row=15; D="0.0006066"; E="-94.26%"; G="19"},
    @{Row=16; D="0.006006"; E="0.57%"; G="19"},
    @{Row=17; E="0.55%"; G="19"},
    @{Row=18; D="3.200"; E="-1.84%"; G="19"},
    @{Row=19; D="2.204"; E="-0.28%"; G="19"},
    @{Row=20; D="0.3106"; E="-1.39%"; G="19"},
    @{Row=21; E="-1.69%"; G="19"},
    @{Row=22; D="3.846"; E="8.83%"; G="19"},
    @{Row=23; D="0.04201"; E="0.71%"; G="19"},
    @{Row=24; D="0.001221"; E="0.17%"; G="19"},
    @{Row=25; D="0.004296"; E="-5.70%"; G="19"},
    @{Row=26; E="-0.11%"; G="19"},
    @{Row=27; D="0.0001509"; E="-22.18%"; G="19"},
    @{Row=28; G="19"},
    @{Row=29; G="19"},
    @{Row=30; G="19"},
    @{Row=31; G="19"},
    @{Row=32; G="19"},
    @{Row=33; G="19"},
    @{Row=34; G="19"},
    @{Row=35; G="19"},
    @{Row=36; G="19"},
    @{Row=37; G="19"},
    @{Row=38; G="19"},
    @{Row=39; G="19"},
    @{Row=40; D="0.03833"; E="-0.01%"; G="19"},
    @{Row=41; B="KickToken"; C="https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"; D="0.006188"; E="11.81%"; G="19"},
    @{Row=42; B="BKEXToken"; C="https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; D="0.1100"; E="-0.79%"; G="19"},
    @{Row=43; D="0.001899"; E="-16.87%"; G="19"},
    @{Row=44; D="0.01179"; E="8.11%"; G="19"},
    @{Row=45; D="0.00005435"; E="-0.08%"; G="19"},
    @{Row=46; E="-0.10%"; G="19"},
    @{Row=47; D="0.06296"; E="-42.34%"; G="19"},
    @{Row=48; D="0.1353"; E="6,215.80%"; G="19"},
    @{Row=49; E="-0.10%"; G="19"},
    @{Row=50; E="-0.10%"; G="19"},
    @{Row=51; G="19"}
)

foreach ($u in $updates) {
    foreach ($col in $u.Keys) {
        if ($col -eq "Row") { continue }
        $cell = $ws.Cells.Item($u.Row, $colIndex[$col])
        if ($col -eq "D" -or $col -eq "E" -or $col -eq "G") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u[$col]
    }
}

Write-Host "Symbol list updated"
